$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (columns A:J were "<name>_old", columns L:U were "<name>_new") ---
# A1:J1 -> suffix "_old" becomes "_FV2410"
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"

# K1 stays "diff" (unchanged)

# L1:U1 -> suffix "_new" becomes "_FV2504"
$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# --- Turn the used range into an Excel Table ---
$tblRange = $ws.Range("A1:U94")
$tbl = $ws.ListObjects.Add(1, $tblRange, $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

Write-Host "done"
